# MotionCoordinates.pptx edit script
# - Update the "datetimeFigureOut" date field text on the slide master and
#   every slide layout from 13/02/2013 -> 01/07/2013.
# - On slide 1:
#     * rotate/re-position the "Straight Arrow Connector 44" connector
#     * delete a handful of now-unused connectors/textboxes
#     * re-position "Rectangle 106"
#     * tidy up run-splitting in the "coordinate system" caption textbox

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder field text: slide master + all slide layouts
# ---------------------------------------------------------------------
$oldDate = "13/02/2013"
$newDate = "01/07/2013"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $msh = $master.Shapes.Item($i)
    if ($msh.HasTextFrame -and $msh.TextFrame.HasText) {
        if ($msh.TextFrame.TextRange.Text -eq $oldDate) {
            $msh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lsh = $layout.Shapes.Item($si)
        if ($lsh.HasTextFrame -and $lsh.TextFrame.HasText) {
            if ($lsh.TextFrame.TextRange.Text -eq $oldDate) {
                $lsh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 1 shape edits
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# 2a. "Straight Arrow Connector 44": flip to the opposite end (rot 180deg)
#     and nudge its position slightly.
#     (Left/Top literals below are the exact float32 values that round-trip
#     to the target EMU offsets of 5476878 / 2704698.)
$conn44 = $s.Shapes.Item("Straight Arrow Connector 44")
$conn44.Rotation = 180
$conn44.Left = 431.250244140625
$conn44.Top = 212.96835327148438

# 2b. Remove shapes that are no longer needed.
$toDelete = @(
    "Straight Arrow Connector 46",
    "Straight Arrow Connector 48",
    "TextBox 52",
    "TextBox 53",
    "Straight Arrow Connector 99",
    "TextBox 100"
)
foreach ($name in $toDelete) {
    $shape = $s.Shapes.Item($name)
    $shape.Delete()
}

# 2c. Reposition "Rectangle 106".
#     (Left/Top literals below are the exact float32 values that round-trip
#     to the target EMU offsets of 5529413 / 2521364.)
$rect106 = $s.Shapes.Item("Rectangle 106")
$rect106.Left = 435.3868713378906
$rect106.Top = 198.53260803222656

# 2d. Merge the "inferior, " / "posterior to anterior, left to right)   "
#     runs (identical formatting) into a single run in the caption textbox.
$caption = $s.Shapes.Item("TextBox 109")
$tr = $caption.TextFrame.TextRange
$paragraphs = $tr.Paragraphs(0, -1)
$secondParagraph = $paragraphs.Paragraphs(2, 1)
$paraText = $secondParagraph.Text
$markerIndex = $paraText.IndexOf("inferior, ")
if ($markerIndex -ge 0) {
    $mergedLength = $paraText.Length - $markerIndex
    $mergeRange = $secondParagraph.Characters($markerIndex + 1, $mergedLength)
    $mergeRange.Text = "inferior, posterior to anterior, left to right)   "
}
